$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/converted-quantity"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
# Row 2 is the "Extension" element; column AI is "Constraint(s)".
# Its combined ele-1/ext-1 constraint text now only applies to the
# "Extension.extension" row (row 4), so the root row's cell is cleared.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
